$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 5-7 (old data no longer needed)
$ws.Rows.Item(5).Resize(3).Delete() | Out-Null

# Row 2: MuSCs / Il10 / Il10rb / ECs
$ws.Range("A2").Value = "MuSCs"
$ws.Range("B2").Value = "Il10"
$ws.Range("C2").Value = "Il10rb"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.08735766666666667
$ws.Range("H2").Value = 0.262073
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 19.741419
$ws.Range("N2").Value = 59.224257
$ws.Range("O2").Value = 0.5456357702458839
$ws.Range("P2").Value = 0.5456357702458838
$ws.Range("Q2").Value = 1.724564300529
$ws.Range("R2").Value = 15.521078704761
$ws.Range("S2").Value = 0.5456357702458839
$ws.Range("T2").Value = 0.5456357702458838

# Row 3: MuSCs / Il10 / Il10rb / FAPs
$ws.Range("A3").Value = "MuSCs"
$ws.Range("B3").Value = "Il10"
$ws.Range("C3").Value = "Il10rb"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.08735766666666667
$ws.Range("H3").Value = 0.262073
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 12.88200833333333
$ws.Range("N3").Value = 38.646025
$ws.Range("O3").Value = 0.3560475839792585
$ws.Range("P3").Value = 0.3560475839792585
$ws.Range("Q3").Value = 1.125342189980556
$ws.Range("R3").Value = 10.128079709825
$ws.Range("S3").Value = 0.3560475839792585
$ws.Range("T3").Value = 0.3560475839792585

# Row 4: MuSCs / Il10 / Il10rb / MuSCs
$ws.Range("A4").Value = "MuSCs"
$ws.Range("B4").Value = "Il10"
$ws.Range("C4").Value = "Il10rb"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.08735766666666667
$ws.Range("H4").Value = 0.262073
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 3.557153333333333
$ws.Range("N4").Value = 10.67146
$ws.Range("O4").Value = 0.09831664577485776
$ws.Range("P4").Value = 0.09831664577485774
$ws.Range("Q4").Value = 0.3107446151755556
$ws.Range("R4").Value = 2.79670153658
$ws.Range("S4").Value = 0.09831664577485776
$ws.Range("T4").Value = 0.09831664577485774
